$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 71; existing rows 71-82 shift down to 72-83.
$ws.Rows(71).Insert()

# Populate the new weekly record in row 71.
$ws.Cells.Item(71, 1).Value = 1
$ws.Cells.Item(71, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(71, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(71, 4).Value = 44543
$ws.Cells.Item(71, 5).Value = 15
$ws.Cells.Item(71, 6).Value = 100112036
$ws.Cells.Item(71, 7).Value = "Caigua"
$ws.Cells.Item(71, 8).Value = "Sin especificar"
$ws.Cells.Item(71, 9).Value = "Primera"
$ws.Cells.Item(71, 10).Value = 120
$ws.Cells.Item(71, 11).Value = 5000
$ws.Cells.Item(71, 12).Value = 6000
$ws.Cells.Item(71, 13).Value = 5500
$ws.Cells.Item(71, 14).Value = "$/caja 20 kilos"
$ws.Cells.Item(71, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(71, 16).Value = 275
$ws.Cells.Item(71, 17).Value = 20
$ws.Cells.Item(71, 18).Value = "Hortaliza"
